$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Copy()
$ws.Range("A41").PasteSpecial(-4122)
$excel.CutCopyMode = $false
